$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume-change columns).
# Price cells (D) are prefixed with a leading apostrophe so Excel's
# smart-entry parsing keeps them as literal text (preserving trailing
# zeros / exact formatting like "142.30") instead of auto-converting
# them to Double values.
$ws.Range("D2").Value = "'25.895.04"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "'1.634.96"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'214.86"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "'0.0631"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "'19.66"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "'1.657.41"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'1.862.18"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "'0.0₃0760"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'25.895.67"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'4.45"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'191.69"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'6.34"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").Value = "'142.30"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").Value = "'1.147.76"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "'0.544"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "'100.81"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'1.772.23"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'0.0₆0111"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'55.60"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +5.86%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'7.63"
$ws.Range("E51").Value = "  +1.32%  "
